$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-range banner) ---
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Row 15 (Rape) ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("M15").Value = 13.888888888888
$ws.Range("N15").Value = -37.878787878787

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = 23.076923076923
$ws.Range("F16").Value = 49
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 44.117647058823
$ws.Range("I16").Value = 401
$ws.Range("J16").Value = 309
$ws.Range("K16").Value = 29.773462783171
$ws.Range("L16").Value = 54.230769230769
$ws.Range("M16").Value = 3.617571059431
$ws.Range("N16").Value = -67.184942716857

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = -3.636363636363
$ws.Range("I17").Value = 632
$ws.Range("J17").Value = 632
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5.158069883527
$ws.Range("M17").Value = 65.879265091863
$ws.Range("N17").Value = -17.493472584856

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 38.095238095238
$ws.Range("I18").Value = 260
$ws.Range("J18").Value = 215
$ws.Range("K18").Value = 20.930232558139
$ws.Range("L18").Value = 24.401913875598
$ws.Range("M18").Value = -18.238993710691
$ws.Range("N18").Value = -84.097859327217

# --- Row 19 (Gr. Larceny) ---
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -11.111111111111
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = 12.068965517241
$ws.Range("I19").Value = 672
$ws.Range("J19").Value = 461
$ws.Range("K19").Value = 45.770065075921
$ws.Range("L19").Value = 76.842105263157
$ws.Range("M19").Value = 146.153846153846
$ws.Range("N19").Value = 39.708939708939

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 63
$ws.Range("H20").Value = -47.619047619047
$ws.Range("I20").Value = 382
$ws.Range("J20").Value = 465
$ws.Range("K20").Value = -17.849462365591
$ws.Range("L20").Value = 66.812227074235
$ws.Range("M20").Value = 31.271477663230
$ws.Range("N20").Value = -72.850035536602

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = -6.153846153846
$ws.Range("F21").Value = 233
$ws.Range("G21").Value = 236
$ws.Range("H21").Value = -1.271186440677
$ws.Range("I21").Value = 2400
$ws.Range("J21").Value = 2132
$ws.Range("K21").Value = 12.570356472795
$ws.Range("L21").Value = 39.049826187717
$ws.Range("M21").Value = 40.597539543058
$ws.Range("N21").Value = -57.173447537473

# --- Row 22 (Transit) ---
$ws.Range("D22").Value = 4
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -87.5
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 75

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "'***.*"
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("L23").Value = 6.593406593406
$ws.Range("M23").Value = 46.969696969697

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 17
$ws.Range("E24").Value = -41.379310344827
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = -16.504854368932
$ws.Range("I24").Value = 1278
$ws.Range("J24").Value = 965
$ws.Range("K24").Value = 32.435233160621
$ws.Range("L24").Value = 48.259860788863
$ws.Range("M24").Value = 89.614243323442

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = -16.901408450704
$ws.Range("I25").Value = 777
$ws.Range("J25").Value = 809
$ws.Range("K25").Value = -3.955500618046
$ws.Range("L25").Value = -13.474387527839
$ws.Range("M25").Value = -8.909730363423

# --- Row 26 (UCR Rape*) ---
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "'***.*"
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = 37.254901960784

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 64
$ws.Range("J27").Value = 72
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 6.666666666666

# --- Row 28 (Shooting Vic.) ---
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = -38.461538461538
$ws.Range("L28").Value = -6.976744186046
$ws.Range("N28").Value = -66.386554621848

# --- Row 29 (Shooting Inc.) ---
$ws.Range("C29").Value = "'0"
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 55
$ws.Range("K29").Value = -32.727272727272
$ws.Range("L29").Value = -5.128205128205
$ws.Range("N29").Value = -66.964285714285

# --- Restore style 14 (general/text, same as column-A label cells) on cells that were
#     just converted to text via the leading quote-prefix above. A plain .Value="'txt"
#     assignment marks the cell as quoted text but also reformats its numFmt/quotePrefix,
#     so we copy the format (only) from the neighboring label cell in column A back onto
#     them to land on the exact same style the workbook already uses for text cells.
$ws.Range("A15").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Range("A23").Copy()
$ws.Range("C23:E23").PasteSpecial(-4122)
$ws.Range("A26").Copy()
$ws.Range("C26:E26").PasteSpecial(-4122)
$ws.Range("A27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("A28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("A29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
